$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 1627.2
$ws.Range("I49").Value = 231
$ws.Range("J49").Value = 2225.5715
$ws.Range("K49").Value = 693
$ws.Range("L49").Value = 6676.7145
$ws.Range("M49").Value = -557
$ws.Range("N49").Value = -6948.7145
$ws.Range("H92").Value = 5354.273
$ws.Range("I92").Value = 4414.7
$ws.Range("J92").Value = 14750
$ws.Range("K92").Value = 4414.7
$ws.Range("L92").Value = 14750
$ws.Range("M92").Value = -3166.7
$ws.Range("N92").Value = -17246
$ws.Range("H107").Value = 1936.5
$ws.Range("I107").Value = 1749.5
$ws.Range("J107").Value = 2497.5
$ws.Range("K107").Value = 1749.5
$ws.Range("L107").Value = 2497.5
$ws.Range("M107").Value = 170.5
$ws.Range("N107").Value = -6337.5
$ws.Range("H116").Value = 6645
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4388616
$ws.Range("I32").Value = 4718791
$ws.Range("K32").Value = 4718791
$ws.Range("M32").Value = -4718504
$ws.Range("H61").Value = 2886682.2
$ws.Range("I61").Value = 3700
$ws.Range("K61").Value = 3700
$ws.Range("M61").Value = -3488
$ws.Range("H74").Value = 14685.685
$ws.Range("I74").Value = 1536.5
$ws.Range("K74").Value = 1536.5
$ws.Range("M74").Value = -662.5
$ws.Range("H77").Value = 14685.685
$ws.Range("I77").Value = 1536.5
$ws.Range("K77").Value = 7682.5
$ws.Range("M77").Value = -3314.5
$ws.Range("H102").Value = 1840.6
$ws.Range("I102").Value = 1426.7273
$ws.Range("J102").Value = 2978.75
$ws.Range("K102").Value = 1426.7273
$ws.Range("L102").Value = 2978.75
$ws.Range("M102").Value = 195.2727
$ws.Range("N102").Value = -6222.75
$ws.Range("H122").Value = 865696.4
$ws.Range("I122").Value = 1152780.5
$ws.Range("K122").Value = 3458341.5
$ws.Range("M122").Value = -3455891.5
$ws.Range("H136").Value = 2886682.2
$ws.Range("I136").Value = 3700
$ws.Range("K136").Value = 11100
$ws.Range("M136").Value = -8550

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66668590
$ws.Range("I86").Value = 1911.5454
$ws.Range("J86").Value = 250001950
$ws.Range("K86").Value = 1911.5454
$ws.Range("L86").Value = 250001950
$ws.Range("M86").Value = -788.5454
$ws.Range("N86").Value = -250004196
$ws.Range("H89").Value = 66668590
$ws.Range("I89").Value = 1911.5454
$ws.Range("J89").Value = 250001950
$ws.Range("K89").Value = 9557.726999999999
$ws.Range("L89").Value = 1250009750
$ws.Range("M89").Value = -3941.726999999999
$ws.Range("N89").Value = -1250020982
$ws.Range("H94").Value = 2856.5
$ws.Range("I94").Value = 2162.818
$ws.Range("K94").Value = 2162.818
$ws.Range("M94").Value = -1711.818

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 42894236
$ws.Range("I132").Value = 2270.318
$ws.Range("K132").Value = 6810.954000000001
$ws.Range("M132").Value = -4280.954000000001
$ws.Range("H134").Value = 19611926
$ws.Range("I134").Value = 967.44116
$ws.Range("K134").Value = 2902.32348
$ws.Range("M134").Value = -367.32348

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 164.71428
$ws.Range("I7").Value = 128.08333
$ws.Range("K7").Value = 384.24999
$ws.Range("M7").Value = -272.24999
$ws.Range("H9").Value = 58656.715
$ws.Range("I9").Value = 101999.75
$ws.Range("J9").Value = 866
$ws.Range("K9").Value = 305999.25
$ws.Range("L9").Value = 2598
$ws.Range("M9").Value = -305775.25
$ws.Range("N9").Value = -3046
$ws.Range("H68").Value = 984.7143
$ws.Range("J68").Value = 1293.6666
$ws.Range("L68").Value = 3880.9998
$ws.Range("N68").Value = -5502.9998
$ws.Range("H71").Value = 984.7143
$ws.Range("J71").Value = 1293.6666
$ws.Range("L71").Value = 11642.9994
$ws.Range("N71").Value = -19754.9994
$ws.Range("H121").Value = 468.54544
$ws.Range("I121").Value = 417.1111
$ws.Range("K121").Value = 1251.3333
$ws.Range("M121").Value = 58.66669999999999
$ws.Range("H131").Value = 1438.95
$ws.Range("I131").Value = 856.3
$ws.Range("J131").Value = 1503.6888
$ws.Range("K131").Value = 2568.9
$ws.Range("L131").Value = 4511.0664
$ws.Range("M131").Value = 2471.1
$ws.Range("N131").Value = -14591.0664

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 122.86207
$ws.Range("I2").Value = 97.05556
$ws.Range("K2").Value = 97.05556
$ws.Range("M2").Value = 15.94444
$ws.Range("H70").Value = 6317.9414
$ws.Range("I70").Value = 4591.4165
$ws.Range("K70").Value = 4591.4165
$ws.Range("M70").Value = -4321.4165
$ws.Range("H73").Value = 6317.9414
$ws.Range("I73").Value = 4591.4165
$ws.Range("K73").Value = 4591.4165
$ws.Range("M73").Value = -3655.4165
$ws.Range("H102").Value = 6145369
$ws.Range("I102").Value = 6759760
$ws.Range("J102").Value = 1457
$ws.Range("K102").Value = 6759760
$ws.Range("L102").Value = 1457
$ws.Range("M102").Value = -6758138
$ws.Range("N102").Value = -4701
$ws.Range("H113").Value = 1413.6471
$ws.Range("J113").Value = 1607.091
$ws.Range("L113").Value = 1607.091
$ws.Range("N113").Value = -5947.091
$ws.Range("H122").Value = 1214342.2
$ws.Range("I122").Value = 1478099.4
$ws.Range("K122").Value = 4434298.199999999
$ws.Range("M122").Value = -4431848.199999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 333334400
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 333334400
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H40").Value = 1156853.9
$ws.Range("I40").Value = 2518.4827
$ws.Range("K40").Value = 2518.4827
$ws.Range("M40").Value = -2382.4827
$ws.Range("H61").Value = 2261.4614
$ws.Range("I61").Value = 1482.1666
$ws.Range("K61").Value = 1482.1666
$ws.Range("M61").Value = -1280.1666
$ws.Range("H74").Value = 40130.668
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 40130.668
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H108").Value = 76325.2
$ws.Range("J108").Value = 76325.2
$ws.Range("L108").Value = 76325.2
$ws.Range("N108").Value = -84005.2
$ws.Range("H113").Value = 2261.4614
$ws.Range("I113").Value = 1482.1666
$ws.Range("K113").Value = 1482.1666
$ws.Range("M113").Value = 687.8334
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9650
$ws.Range("I62").Value = 12900.333
$ws.Range("K62").Value = 12900.333
$ws.Range("M62").Value = -12276.333
$ws.Range("H65").Value = 9650
$ws.Range("I65").Value = 12900.333
$ws.Range("K65").Value = 64501.665
$ws.Range("M65").Value = -61381.665
$ws.Range("H113").Value = 2298.5151
$ws.Range("I113").Value = 2253.1614
$ws.Range("J113").Value = 3001.5
$ws.Range("K113").Value = 6759.4842
$ws.Range("L113").Value = 9004.5
$ws.Range("M113").Value = -4589.4842
$ws.Range("N113").Value = -13344.5
